{"js": "// Apply the \"Gradient + p\u0159ep\u00edn\u00e1n\u00ed gradientu\" edit:\n//  - row \"interpolace barvy ur\u010den\u00e9 ve vrcholech\" / \"0,5\" / \"(Ve 2d)\":\n//      * points value \"0,5\" -> \"1\"\n//      * note text \"(Ve 2d)\" -> \"Tla\u010d\u00edtko V na jehlanu\"\n//      * row height 300 twips (15pt) -> 395 twips (19.75pt)\n\n// --- 1) \"0,5\" -> \"1\" -------------------------------------------------\nconst pointsResults = context.document.body.search(\"0,5\", { matchCase: true, matchWholeWord: false });\npointsResults.load(\"items,text\");\nawait context.sync();\n\nfor (const item of pointsResults.items) {\n  item.insertText(\"1\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) \"(Ve 2d)\" -> \"Tla\u010d\u00edtko V na jehlanu\" --------------------------\nconst noteResults = context.document.body.search(\"(Ve 2d)\", { matchCase: true, matchWholeWord: false });\nnoteResults.load(\"items,text\");\nawait context.sync();\n\nfor (const item of noteResults.items) {\n  item.insertText(\"Tla\u010d\u00edtko V na jehlanu\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 3) Grow the row's height from 300 twips (15pt) to 395 twips (19.75pt)\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.load(\"values\");\n}\nawait context.sync();\n\nconst targetRow = rows.items.find((row) =>\n  row.values.some((rowValues) => rowValues.some((cellText) => cellText.indexOf(\"interpolace barvy\") !== -1))\n);\n\nif (targetRow) {\n  targetRow.preferredHeight = 19.75; // 395 twips == 19.75 points\n}\nawait context.sync();\n", "ps1": "# Apply the \"Gradient + p\u0159ep\u00edn\u00e1n\u00ed gradientu\" edit:\n#  - row \"interpolace barvy ur\u010den\u00e9 ve vrcholech\" / \"0,5\" / \"(Ve 2d)\":\n#      * points value \"0,5\" -> \"1\"\n#      * note text \"(Ve 2d)\" -> \"Tla\u010d\u00edtko V na jehlanu\"\n#      * row height 300 twips (15pt) -> 395 twips (19.75pt)\n\n$d = $word.ActiveDocument\n\n# --- 1) \"0,5\" -> \"1\" ---------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"0,5\"\n$find.Replacement.Text = \"1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# --- 2) \"(Ve 2d)\" -> \"Tla\u010d\u00edtko V na jehlanu\" ----------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"(Ve 2d)\"\n$find2.Replacement.Text = \"Tla\u010d\u00edtko V na jehlanu\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# --- 3) Grow the row's height from 300 twips (15pt) to 395 twips (19.75pt)\n$table = $d.Tables.Item(1)\nforeach ($row in $table.Rows) {\n    if ($row.Cells.Item(1).Range.Text -like \"*interpolace barvy*\") {\n        $row.Height = 19.75   # 395 twips == 19.75 points\n        break\n    }\n}\n"}
